$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 528.9167
$ws.Range("I39").Value = 125
$ws.Range("J39").Value = 730.875
$ws.Range("K39").Value = 375
$ws.Range("L39").Value = 2192.625
$ws.Range("M39").Value = -79
$ws.Range("N39").Value = -2784.625

$ws.Range("H127").Value = 1886.8064
$ws.Range("I127").Value = 850.6667
$ws.Range("K127").Value = 2552.0001
$ws.Range("M127").Value = 2407.9999

$ws.Range("H129").Value = 1174.8334
$ws.Range("J129").Value = 1705.4286
$ws.Range("L129").Value = 5116.2858
$ws.Range("N129").Value = -15116.2858

$ws.Range("H138").Value = 6537
$ws.Range("I138").Value = 1181.5385
$ws.Range("J138").Value = 9437.875
$ws.Range("K138").Value = 3544.6155
$ws.Range("L138").Value = 28313.625
$ws.Range("M138").Value = 1595.3845
$ws.Range("N138").Value = -38593.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5630.772
$ws.Range("I32").Value = 4207.2793
$ws.Range("J32").Value = 10002.929
$ws.Range("K32").Value = 4207.2793
$ws.Range("L32").Value = 10002.929
$ws.Range("M32").Value = -3920.2793
$ws.Range("N32").Value = -10576.929

$ws.Range("H76").Value = 30057.6
$ws.Range("J76").Value = 30057.6
$ws.Range("L76").Value = 30057.6
$ws.Range("N76").Value = -30733.6

$ws.Range("H79").Value = 30057.6
$ws.Range("J79").Value = 30057.6
$ws.Range("L79").Value = 30057.6
$ws.Range("N79").Value = -32397.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 13893347
$ws.Range("I99").Value = 1253
$ws.Range("J99").Value = 25007022
$ws.Range("K99").Value = 1253
$ws.Range("L99").Value = 25007022
$ws.Range("M99").Value = 245
$ws.Range("N99").Value = -25010018

$ws.Range("H122").Value = 833
$ws.Range("J122").Value = 582.5
$ws.Range("L122").Value = 1747.5
$ws.Range("N122").Value = -6647.5

$ws.Range("H126").Value = 13893347
$ws.Range("I126").Value = 1253
$ws.Range("J126").Value = 25007022
$ws.Range("K126").Value = 3759
$ws.Range("L126").Value = 75021066
$ws.Range("M126").Value = -1289
$ws.Range("N126").Value = -75026006

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 316659.16
$ws.Range("I5").Value = 578
$ws.Range("J5").Value = 546536.4
$ws.Range("K5").Value = 1734
$ws.Range("L5").Value = 1639609.2
$ws.Range("M5").Value = -1622
$ws.Range("N5").Value = -1639833.2

$ws.Range("H122").Value = 699.9474
$ws.Range("I122").Value = 561.25
$ws.Range("K122").Value = 5051.25
$ws.Range("M122").Value = -2601.25

$ws.Range("H129").Value = 2202
$ws.Range("I129").Value = 1503.625
$ws.Range("J129").Value = 3319.4
$ws.Range("K129").Value = 4510.875
$ws.Range("L129").Value = 9958.2
$ws.Range("M129").Value = 489.125
$ws.Range("N129").Value = -19958.2

$ws.Range("H135").Value = 316659.16
$ws.Range("I135").Value = 578
$ws.Range("J135").Value = 546536.4
$ws.Range("K135").Value = 5202
$ws.Range("L135").Value = 4918827.600000001
$ws.Range("M135").Value = -2667
$ws.Range("N135").Value = -4923897.600000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6299.0356
$ws.Range("I70").Value = 6342.304
$ws.Range("J70").Value = 6100
$ws.Range("K70").Value = 6342.304
$ws.Range("L70").Value = 6100
$ws.Range("M70").Value = -6072.304
$ws.Range("N70").Value = -6640

$ws.Range("H73").Value = 6299.0356
$ws.Range("I73").Value = 6342.304
$ws.Range("J73").Value = 6100
$ws.Range("K73").Value = 6342.304
$ws.Range("L73").Value = 6100
$ws.Range("M73").Value = -5406.304
$ws.Range("N73").Value = -7972

$ws.Range("H80").Value = 2521.0688
$ws.Range("I80").Value = 2448.2
$ws.Range("J80").Value = 2976.5
$ws.Range("K80").Value = 2448.2
$ws.Range("L80").Value = 2976.5
$ws.Range("M80").Value = -1450.2
$ws.Range("N80").Value = -4972.5

$ws.Range("H83").Value = 2521.0688
$ws.Range("I83").Value = 2448.2
$ws.Range("J83").Value = 2976.5
$ws.Range("K83").Value = 12241
$ws.Range("L83").Value = 14882.5
$ws.Range("M83").Value = -7249
$ws.Range("N83").Value = -24866.5

$ws.Range("H102").Value = 1358.4
$ws.Range("I102").Value = 1341.2
$ws.Range("J102").Value = 1392.8
$ws.Range("K102").Value = 1341.2
$ws.Range("L102").Value = 1392.8
$ws.Range("M102").Value = 280.8
$ws.Range("N102").Value = -4636.8

$ws.Range("H103").Value = 38414.57
$ws.Range("J103").Value = 38414.57
$ws.Range("L103").Value = 38414.57
$ws.Range("N103").Value = -40758.57

$ws.Range("H136").Value = 38000
$ws.Range("J136").Value = 38000
$ws.Range("L136").Value = 114000
$ws.Range("N136").Value = -119100

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 9803
$ws.Range("I33").Value = 9015
$ws.Range("K33").Value = 9015
$ws.Range("M33").Value = -8725

$ws.Range("H43").Value = 5500
$ws.Range("J43").Value = 5500
$ws.Range("L43").Value = 5500
$ws.Range("N43").Value = -5886

$ws.Range("H61").Value = 2537.6316
$ws.Range("I61").Value = 2185
$ws.Range("J61").Value = 3142.1428
$ws.Range("K61").Value = 2185
$ws.Range("L61").Value = 3142.1428
$ws.Range("M61").Value = -1983
$ws.Range("N61").Value = -3546.1428

$ws.Range("H64").Value = 27999
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 27999
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H68").Value = 76925200
$ws.Range("I68").Value = 1763.6
$ws.Range("J68").Value = 125002344
$ws.Range("K68").Value = 1763.6
$ws.Range("L68").Value = 125002344
$ws.Range("M68").Value = -1014.6
$ws.Range("N68").Value = -125003842

$ws.Range("H71").Value = 76925200
$ws.Range("I71").Value = 1763.6
$ws.Range("J71").Value = 125002344
$ws.Range("K71").Value = 8818
$ws.Range("L71").Value = 625011720
$ws.Range("M71").Value = -5074
$ws.Range("N71").Value = -625019208

$ws.Range("H113").Value = 2537.6316
$ws.Range("I113").Value = 2185
$ws.Range("J113").Value = 3142.1428
$ws.Range("K113").Value = 2185
$ws.Range("L113").Value = 3142.1428
$ws.Range("M113").Value = -15
$ws.Range("N113").Value = -7482.1428

$ws.Range("H132").Value = 16056277
$ws.Range("I132").Value = 24082040
$ws.Range("J132").Value = 4749.222
$ws.Range("K132").Value = 72246120
$ws.Range("L132").Value = 14247.666
$ws.Range("M132").Value = -72243590
$ws.Range("N132").Value = -19307.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1632.1428
$ws.Range("I113").Value = 1590.3
$ws.Range("J113").Value = 1736.75
$ws.Range("K113").Value = 4770.9
$ws.Range("L113").Value = 5210.25
$ws.Range("M113").Value = -2600.9
$ws.Range("N113").Value = -9550.25

$ws.Range("H122").Value = 2426
$ws.Range("I122").Value = 2502
$ws.Range("J122").Value = 2350
$ws.Range("K122").Value = 7506
$ws.Range("L122").Value = 7050
$ws.Range("M122").Value = -5056
$ws.Range("N122").Value = -11950

$ws.Range("H126").Value = 1382.3334
$ws.Range("I126").Value = 948.7143
$ws.Range("J126").Value = 2900
$ws.Range("K126").Value = 2846.1429
$ws.Range("L126").Value = 8700
$ws.Range("M126").Value = -376.1428999999998
$ws.Range("N126").Value = -13640

$ws.Range("H128").Value = 47857.855
$ws.Range("J128").Value = 47857.855
$ws.Range("L128").Value = 47857.855
$ws.Range("N128").Value = -57817.855
